# RP_titul_list.docx - apply the edits described by the commit
# "Úprava dokumentace + návrh databáze"

$d = $word.ActiveDocument

# 1) Téma (topic) line: "Úbytek konipasů v Plzeňském kraji" ->
#    "Databáze seriálů (placeholder název)"
$d.Content.Find.Execute(
    "Úbytek konipasů v Plzeňském kraji", $true, $false, $false, $false, $false,
    $true, 1, $false, "Databáze seriálů (placeholder název)", 2) | Out-Null

# 2) Autor práce (author): "Valerián DOBROTIVÝ" ->
#    "Adam Vlček, Ondřej Pták, Jan Rehák"
$d.Content.Find.Execute(
    "Valerián DOBROTIVÝ", $true, $false, $false, $false, $false,
    $true, 1, $false, "Adam Vlček, Ondřej Pták, Jan Rehák", 2) | Out-Null

# 3) Obor studia (field of study): collapse the long "78-42-M/01 Technické
#    lyceum (nebo 26-41-M/01 Elektrotechnika nebo 18-20-M/01  Inf.
#    technologie)" text down to "18-20-M/01 Inf. technologie"
$d.Content.Find.Execute(
    "78-42-M/01 Technické lyceum (nebo 26-41-M/01 Elektrotechnika nebo 18-20-M/01  Inf. technologie)",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "18-20-M/01 Inf. technologie", 2) | Out-Null

# 4) Třída (class): "4. X" -> "2. H"
$d.Content.Find.Execute(
    "4. X", $true, $false, $false, $false, $false,
    $true, 1, $false, "2. H", 2) | Out-Null

# 5) Předmět (subject): "Matematika" -> "Vývoj webových stránek"
$d.Content.Find.Execute(
    "Matematika", $true, $false, $false, $false, $false,
    $true, 1, $false, "Vývoj webových stránek", 2) | Out-Null

# 6) Zadávající učitel (assigning teacher): "Ing. Alfréd KONIPÁSEK" ->
#    "Mgr. Pavlína Lukešová"
$d.Content.Find.Execute(
    "Ing. Alfréd KONIPÁSEK", $true, $false, $false, $false, $false,
    $true, 1, $false, "Mgr. Pavlína Lukešová", 2) | Out-Null

# 7) Dne (date): "30. 3. 2020" -> "XX. X. 2023" (also removes the
#    now-orphaned _GoBack bookmark that used to sit inside this text)
$d.Content.Find.Execute(
    "30. 3. 2020", $true, $false, $false, $false, $false,
    $true, 1, $false, "XX. X. 2023", 2) | Out-Null
